# Updated cryptos list on Fri Dec  1 04:31:03 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells whose new value parses as a plain decimal number
# must be forced to stay text (matching the source data, which stores every
# price as a literal string) -- otherwise Excel's smart-typing would turn
# e.g. "100.00" into the number 100 and drop the trailing zeros.
$textPriceCells = @(
    "D5","D6","D10","D13","D14","D16","D19","D20","D22","D24","D26","D27",
    "D28","D29","D30","D35","D36","D37","D38","D40","D42","D43","D49","D50"
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "38.194.43"
$ws.Range("E2").Value = "  +0.36%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.098.19"
$ws.Range("E3").Value = "  +2.76%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "229.75"
$ws.Range("E5").Value = "  +0.28%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +0.18%  "

# Row 7 - Solana
$ws.Range("E7").Value = "  +0.13%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.75%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0842"
$ws.Range("E10").Value = "  +2.31%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.68%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.405.63"
$ws.Range("E12").Value = "  +2.65%  "

# Row 13 - Avalanche
$ws.Range("D13").Value = "22.46"
$ws.Range("E13").Value = "  +4.91%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "14.67"
$ws.Range("E14").Value = "  -0.36%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  +5.77%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.776"
$ws.Range("E16").Value = "  +0.72%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.095.87"
$ws.Range("E17").Value = "  +1.70%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "38.143.29"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "6.01"

# Row 20 - Litecoin
$ws.Range("D20").Value = "70.25"
$ws.Range("E20").Value = "  +0.33%  "

# Row 21 - ShibaInu
$ws.Range("D21").Value = "0.0₃0837"
$ws.Range("E21").Value = "  +1.03%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "224.36"
$ws.Range("E22").Value = "  -0.29%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.54%  "

# Row 24 - Toncoin
$ws.Range("D24").Value = "2.44"

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +3.08%  "

# Row 26 - Monero
$ws.Range("D26").Value = "170.28"
$ws.Range("E26").Value = "  +1.85%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "9.43"
$ws.Range("E27").Value = "  +0.83%  "

# Row 28 - Kaspa
$ws.Range("D28").Value = "0.132"
$ws.Range("E28").Value = "  +1.77%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "19.06"
$ws.Range("E29").Value = "  +0.44%  "

# Row 30 - ImmutableX
$ws.Range("D30").Value = "1.34"
$ws.Range("E30").Value = "  +4.15%  "

# Row 31 - Stellar
$ws.Range("E31").Value = "  -0.34%  "

# Row 32 - WEMIXToken
$ws.Range("E32").Value = "  +9.52%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +2.86%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +0.04%  "

# Row 35 - Hedera
$ws.Range("D35").Value = "0.0607"
$ws.Range("E35").Value = "  -0.41%  "

# Row 36/37 - LidoDAOToken and THORChain swap rank positions
$ws.Range("B36").Value = "THORChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D36").Value = "6.50"
$ws.Range("E36").Value = "  +0.12%  "

$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.41"
$ws.Range("E37").Value = "  +5.41%  "

# Row 38 - RenderToken
$ws.Range("D38").Value = "3.56"
$ws.Range("E38").Value = "  +7.67%  "

# Row 39 - BinanceUSD
$ws.Range("E39").Value = "  -0.23%  "

# Row 40 - InjectiveProtocol
$ws.Range("D40").Value = "17.99"
$ws.Range("E40").Value = "  +0.50%  "

# Row 41 - Maker
$ws.Range("D41").Value = "1.547.72"
$ws.Range("E41").Value = "  +1.23%  "

# Row 42 - Aave
$ws.Range("D42").Value = "100.00"
$ws.Range("E42").Value = "  +3.51%  "

# Row 43 - VeChain
$ws.Range("D43").Value = "0.0219"
$ws.Range("E43").Value = "  -0.11%  "

# Row 44 - HuobiToken
$ws.Range("E44").Value = "  +0.89%  "

# Row 45 - Cronos
$ws.Range("E45").Value = "  -1.43%  "

# Row 46 - FTXToken
$ws.Range("E46").Value = "  +1.67%  "

# Row 47 - TrustWalletToken
$ws.Range("E47").Value = "  +0.48%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  +1.19%  "

# Row 49/50 - MXToken and FraxShare swap rank positions
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.27"
$ws.Range("E49").Value = "  +1.79%  "

$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "3.01"
$ws.Range("E50").Value = "  +1.27%  "

# Row 51 - RocketPoolETH
$ws.Range("D51").Value = "2.292.99"
$ws.Range("E51").Value = "  +2.71%  "
